$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Closing_Price")

# Update existing values in B196:B198
$ws.Range("B196").Value = 114976911118419.2
$ws.Range("B197").Value = 114672959226172
$ws.Range("B198").Value = 112917914904701.7

# Add new row 199
$ws.Range("A199").Value = 45078
$ws.Range("A199").Style = $ws.Range("A198").Style
$ws.Range("B199").Value = 112569394359339.9
